$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") cell updates per repulled data / mean calculation
$updates = @{
    "F3"  = -1
    "F4"  = 2
    "F13" = 0
    "F19" = -3
    "F20" = -2
    "F21" = -2
    "F25" = 0
    "F28" = 2
    "F29" = -1
    "F36" = -1
    "F37" = 0
    "F43" = -3
    "F47" = 2
    "F50" = 0
    "F61" = 1
    "F62" = -3
    "F64" = -6
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
